$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1081.909
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("M2").Value = 13
$ws.Range("H4").Value = 1623.2
$ws.Range("I4").Value = 1054.5
$ws.Range("K4").Value = 1054.5
$ws.Range("M4").Value = -940.5
$ws.Range("H17").Value = 2364.9033
$ws.Range("J17").Value = 2364.9033
$ws.Range("L17").Value = 7094.7099
$ws.Range("N17").Value = -7430.7099
$ws.Range("H93").Value = 19499.5
$ws.Range("J93").Value = 19499.5
$ws.Range("L93").Value = 19499.5
$ws.Range("N93").Value = -24491.5
$ws.Range("H106").Value = 3290.1
$ws.Range("I106").Value = 2962.625
$ws.Range("K106").Value = 2962.625
$ws.Range("M106").Value = -2331.625
$ws.Range("H107").Value = 620.5833
$ws.Range("I107").Value = 620.5833
$ws.Range("K107").Value = 620.5833
$ws.Range("M107").Value = 1299.4167
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("H132").Value = 2171.1628
$ws.Range("I132").Value = 2141.0789
$ws.Range("J132").Value = 2399.8
$ws.Range("K132").Value = 6423.236699999999
$ws.Range("L132").Value = 7199.400000000001
$ws.Range("M132").Value = -3893.236699999999
$ws.Range("N132").Value = -12259.4
$ws.Range("M118").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 48666.668
$ws.Range("J92").Value = 62000
$ws.Range("L92").Value = 62000
$ws.Range("N92").Value = -66992
$ws.Range("H110").Value = 2728.8
$ws.Range("I110").Value = 2786.4194
$ws.Range("K110").Value = 2786.4194
$ws.Range("M110").Value = -741.4194000000002
$ws.Range("H132").Value = 6887.3887
$ws.Range("I132").Value = 1997.9166
$ws.Range("J132").Value = 16666.334
$ws.Range("K132").Value = 5993.7498
$ws.Range("L132").Value = 49999.00199999999
$ws.Range("M132").Value = -3463.7498
$ws.Range("N132").Value = -55059.00199999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1609.7142
$ws.Range("I22").Value = 1378
$ws.Range("K22").Value = 1378
$ws.Range("M22").Value = -1205
$ws.Range("H94").Value = 733.85364
$ws.Range("I94").Value = 729.4865
$ws.Range("J94").Value = 774.25
$ws.Range("K94").Value = 729.4865
$ws.Range("L94").Value = 774.25
$ws.Range("M94").Value = -278.4865
$ws.Range("N94").Value = -1676.25
$ws.Range("H99").Value = 2152.2307
$ws.Range("I99").Value = 1998.091
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 1998.091
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -500.0909999999999
$ws.Range("N99").Value = -5996
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("H105").Value = 3879.4211
$ws.Range("J105").Value = 3584.75
$ws.Range("L105").Value = 3584.75
$ws.Range("N105").Value = -7078.75
$ws.Range("H106").Value = 34665
$ws.Range("J106").Value = 34665
$ws.Range("L106").Value = 34665
$ws.Range("N106").Value = -37189
$ws.Range("H134").Value = 5578.0625
$ws.Range("I134").Value = 1973.9166
$ws.Range("J134").Value = 16390.5
$ws.Range("K134").Value = 5921.7498
$ws.Range("L134").Value = 49171.5
$ws.Range("M134").Value = -3386.7498
$ws.Range("N134").Value = -54241.5
$ws.Range("M102").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7723.7334
$ws.Range("I31").Value = 2863.0417
$ws.Range("J31").Value = 27166.5
$ws.Range("K31").Value = 2863.0417
$ws.Range("L31").Value = 27166.5
$ws.Range("M31").Value = -2568.0417
$ws.Range("N31").Value = -27756.5
$ws.Range("H34").Value = 7723.7334
$ws.Range("I34").Value = 2863.0417
$ws.Range("J34").Value = 27166.5
$ws.Range("K34").Value = 2863.0417
$ws.Range("L34").Value = 27166.5
$ws.Range("M34").Value = -2661.0417
$ws.Range("N34").Value = -27570.5
$ws.Range("H105").Value = 4544
$ws.Range("I105").Value = 3680
$ws.Range("K105").Value = 3680
$ws.Range("M105").Value = -1933
$ws.Range("H132").Value = 3266.5833
$ws.Range("I132").Value = 3266.5833
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9799.749899999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7269.749899999999
$ws.Range("H134").Value = 6665.3335
$ws.Range("I134").Value = 5152.5386
$ws.Range("K134").Value = 15457.6158
$ws.Range("M134").Value = -12922.6158
$ws.Range("N132").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1221.25
$ws.Range("I5").Value = 380
$ws.Range("J5").Value = 1501.6666
$ws.Range("K5").Value = 1140
$ws.Range("L5").Value = 4504.9998
$ws.Range("M5").Value = -1028
$ws.Range("N5").Value = -4728.9998
$ws.Range("H120").Value = 13281
$ws.Range("I120").Value = 9184.333000000001
$ws.Range("J120").Value = 22498.5
$ws.Range("K120").Value = 27552.999
$ws.Range("L120").Value = 67495.5
$ws.Range("M120").Value = -22714.999
$ws.Range("N120").Value = -77171.5
$ws.Range("H122").Value = 863.7143
$ws.Range("J122").Value = 858
$ws.Range("L122").Value = 7722
$ws.Range("N122").Value = -12622
$ws.Range("H135").Value = 1221.25
$ws.Range("I135").Value = 380
$ws.Range("J135").Value = 1501.6666
$ws.Range("K135").Value = 3420
$ws.Range("L135").Value = 13514.9994
$ws.Range("M135").Value = -885
$ws.Range("N135").Value = -18584.9994

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 31336
$ws.Range("J104").Value = 31336
$ws.Range("L104").Value = 31336
$ws.Range("N104").Value = -38324
$ws.Range("H122").Value = 5430.2
$ws.Range("I122").Value = 995.4
$ws.Range("J122").Value = 14299.8
$ws.Range("K122").Value = 2986.2
$ws.Range("L122").Value = 42899.39999999999
$ws.Range("M122").Value = -536.1999999999998
$ws.Range("N122").Value = -47799.39999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1943.5807
$ws.Range("I16").Value = 1820.25
$ws.Range("J16").Value = 2366.4285
$ws.Range("K16").Value = 1820.25
$ws.Range("L16").Value = 2366.4285
$ws.Range("M16").Value = -1650.25
$ws.Range("N16").Value = -2706.4285
$ws.Range("H40").Value = 5484.6665
$ws.Range("I40").Value = 4346.909
$ws.Range("K40").Value = 4346.909
$ws.Range("M40").Value = -4210.909
$ws.Range("H68").Value = 5195.2856
$ws.Range("I68").Value = 5151.6665
$ws.Range("J68").Value = 5228
$ws.Range("K68").Value = 5151.6665
$ws.Range("L68").Value = 5228
$ws.Range("M68").Value = -4402.6665
$ws.Range("N68").Value = -6726
$ws.Range("H71").Value = 5195.2856
$ws.Range("I71").Value = 5151.6665
$ws.Range("J71").Value = 5228
$ws.Range("K71").Value = 25758.3325
$ws.Range("L71").Value = 26140
$ws.Range("M71").Value = -22014.3325
$ws.Range("N71").Value = -33628
$ws.Range("H122").Value = 4790.64
$ws.Range("I122").Value = 4225.727
$ws.Range("J122").Value = 8933.333000000001
$ws.Range("K122").Value = 12677.181
$ws.Range("L122").Value = 26799.999
$ws.Range("M122").Value = -10227.181
$ws.Range("N122").Value = -31699.999
$ws.Range("H136").Value = 3846.7437
$ws.Range("I136").Value = 1843.9615
$ws.Range("J136").Value = 7852.3076
$ws.Range("K136").Value = 5531.8845
$ws.Range("L136").Value = 23556.9228
$ws.Range("M136").Value = -2981.8845
$ws.Range("N136").Value = -28656.9228

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6965
$ws.Range("I81").Value = 8000
$ws.Range("J81").Value = 6447.5
$ws.Range("K81").Value = 16000
$ws.Range("L81").Value = 12895
$ws.Range("M81").Value = -14939
$ws.Range("N81").Value = -15017
$ws.Range("H84").Value = 6965
$ws.Range("I84").Value = 8000
$ws.Range("J84").Value = 6447.5
$ws.Range("K84").Value = 80000
$ws.Range("L84").Value = 64475
$ws.Range("M84").Value = -74696
$ws.Range("N84").Value = -75083
$ws.Range("H107").Value = 1367.84
$ws.Range("I107").Value = 1372.0625
$ws.Range("J107").Value = 1360.3334
$ws.Range("K107").Value = 4116.1875
$ws.Range("L107").Value = 4081.0002
$ws.Range("M107").Value = -2196.1875
$ws.Range("N107").Value = -7921.0002
$ws.Range("H132").Value = 2866.7585
$ws.Range("I132").Value = 2116.1482
$ws.Range("K132").Value = 6348.444600000001
$ws.Range("M132").Value = -3818.444600000001
